$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201, pushing the existing rows 201:234 down to 202:235
$ws.Rows(201).Insert()

# Fill in the newly inserted row 201 with the new weekly record
$ws.Cells.Item(201, 1).Value  = 7
$ws.Cells.Item(201, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(201, 3).Value  = "Ñuble"
$ws.Cells.Item(201, 4).Value  = 44816
$ws.Cells.Item(201, 5).Value  = 16
$ws.Cells.Item(201, 6).Value  = 100112017
$ws.Cells.Item(201, 7).Value  = "Apio"
$ws.Cells.Item(201, 8).Value  = "Americana (o)"
$ws.Cells.Item(201, 9).Value  = "Primera"
$ws.Cells.Item(201, 10).Value = 120
$ws.Cells.Item(201, 11).Value = 9000
$ws.Cells.Item(201, 12).Value = 10000
$ws.Cells.Item(201, 13).Value = 9500
$ws.Cells.Item(201, 14).Value = "$/docena de matas"
$ws.Cells.Item(201, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(201, 16).Value = 1583
$ws.Cells.Item(201, 17).Value = 6
$ws.Cells.Item(201, 18).Value = "Hortaliza"
